$wb = $excel.ActiveWorkbook

# --- "Deep water port" sheet: row 35 already exists but the BR total
# formulas in D3/E3 were never extended to include it. Fix the ranges.
$wsDeep = $wb.Worksheets.Item("Deep water port")
$wsDeep.Range("D3").Formula = "=SUM(D4:D35)"
$wsDeep.Range("E3").Formula = "=SUM(E4:E35)"

# --- "Shallow water port" sheet: same fix, extend BR total formulas to
# include row 21, and correct several BR values in column C.
$wsShallow = $wb.Worksheets.Item("Shallow water port")
$wsShallow.Range("D3").Formula = "=SUM(D4:D21)"
$wsShallow.Range("E3").Formula = "=SUM(E4:E21)"

$wsShallow.Range("C4").Value = 100   # Brig
$wsShallow.Range("C5").Value = 100   # Rattlesnake Heavy
$wsShallow.Range("C13").Value = 80   # Hercules
$wsShallow.Range("C14").Value = 70   # Pandora
$wsShallow.Range("C15").Value = 55   # Cutter
$wsShallow.Range("C21").Value = 50   # Pickle
